# Apply updates to column F (dSF) for specific rows on the active worksheet,
# as described in the commit: "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    14 = -5
    18 = 0
    23 = 0
    26 = -10
    34 = -2
    36 = -2
    37 = 0
    42 = -5
    43 = -4
    46 = -5
    47 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
